$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the degree abbreviations to include periods
$ws.Range("A16").Value = "M.Sc. Biology"
$ws.Range("A17").Value = "B.Sc. Terrestrial & Aquatic Eology (Honours)"

# Reflect the cursor/selection location left after the edit
$ws.Range("A17").Select()
